$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("peliculas o documentales")
$ws2 = $wb.Worksheets.Item("series")

# Helper: write a "true" numeric value into a cell that otherwise carries a
# text ("@") number format, so it is persisted as a literal <v> number
# instead of being coerced into a shared string.
function Set-NumericValue {
    param($range, $value)
    $range.Style = "Normal"
    $range.Value = $value
    $range.NumberFormat = "@"
}

# ----------------------------------------------------------------------
# Sheet "peliculas o documentales": add two new rows (3 and 4)
# Row 2 already holds "Guerrilla del Oro" and stays untouched.
# ----------------------------------------------------------------------

# Row 3: Planeta Hostil / WILD / Premium / Documental / HD / 2019
$ws1.Range("A3").Value = "Planeta Hostil"
$ws1.Range("B3").Value = "WILD"
$ws1.Range("C3").Value = "Premium"
$ws1.Range("D3").Value = "Documental"
$ws1.Range("E3").Value = "HD"
Set-NumericValue $ws1.Range("F3") 2019

# Row 4: Sumergidos / (no proveedor) / Arriendo / Acción / HD / "2019" / $3.490
$ws1.Range("A4").Value = "Sumergidos"
$ws1.Range("C4").Value = "Arriendo"
$ws1.Range("D4").Value = "Acción"
$ws1.Range("E4").Value = "HD"
$ws1.Range("F4").Value = "2019"
$ws1.Range("G4").Value = "`$3.490"

# ----------------------------------------------------------------------
# Sheet "series": rewrite row 2 (same content, numeric AGNO/EPISODIOS)
# and add new row 3 for "The Outsider - El Visitante"
# ----------------------------------------------------------------------

$ws2.Range("A2").Value = "My Brilliant Friend"
$ws2.Range("C2").Value = "Gratis"
$ws2.Range("D2").Value = "Drama"
$ws2.Range("E2").Value = "HD"
Set-NumericValue $ws2.Range("F2") 2020
$ws2.Range("G2").Value = "Temporada 02"
Set-NumericValue $ws2.Range("H2") 1

$ws2.Range("A3").Value = "The Outsider - El Visitante"
$ws2.Range("C3").Value = "Gratis"
$ws2.Range("D3").Value = "Acción"
$ws2.Range("E3").Value = "HD"
Set-NumericValue $ws2.Range("F3") 2020
$ws2.Range("G3").Value = "Temporada 01"
Set-NumericValue $ws2.Range("H3") 10
